# Add Test Data for Croatia Market
# - Duplicate the "Turkey" sheet (keeping it last, before the new one) to
#   create a new "Croatia" sheet with the same layout/styling.
# - Update the market name / ticket reference cells on the new sheet.
# - Leave the new sheet as the active / selected tab, with B4 selected.

$wb = $excel.ActiveWorkbook

$turkey = $wb.Worksheets.Item("Turkey")

# Select a cell on Turkey before duplicating it, matching the final
# selection state left behind on that sheet.
$turkey.Activate()
$turkey.Range("I24").Select()

# Duplicate Turkey into a new sheet placed right after it.
$turkey.Copy([System.Reflection.Missing]::Value, $turkey)
$croatia = $wb.Worksheets.Item($turkey.Index + 1)
$croatia.Name = "Croatia"

# Update the market name and reference values for Croatia.
$croatia.Range("B2").Value = "Croatia Market"
$croatia.Range("B4").Value = "NGC-3139/T2421"

# Leave the new sheet active with B4 selected.
$croatia.Activate()
$croatia.Range("B4").Select()
